# Update cryptos list cell values to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.074.48"
$ws.Range("E2").Value = "  -0.68%  "

$ws.Range("D3").Value = "'3.299.58"
$ws.Range("E3").Value = "  -1.91%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'188.97"
$ws.Range("E5").Value = "  +3.14%  "

$ws.Range("D6").Value = "'558.57"
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  -1.95%  "

$ws.Range("D9").Value = "'3.289.83"
$ws.Range("E9").Value = "  -1.98%  "

$ws.Range("D10").Value = "'0.183"
$ws.Range("E10").Value = "  -1.93%  "

$ws.Range("D11").Value = "'0.586"
$ws.Range("E11").Value = "  -1.62%  "

$ws.Range("D12").Value = "'47.62"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("D13").Value = "'0.0000269"
$ws.Range("E13").Value = "  +0.33%  "

$ws.Range("D14").Value = "'8.65"
$ws.Range("E14").Value = "  -1.09%  "

$ws.Range("D15").Value = "'3.834.69"
$ws.Range("E15").Value = "  -1.80%  "

$ws.Range("D16").Value = "'605.30"
$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").Value = "'66.117.18"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").Value = "'17.98"
$ws.Range("E18").Value = "  -1.48%  "

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").Value = "'3.308.96"
$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").Value = "'11.05"
$ws.Range("E21").Value = "  -4.26%  "

$ws.Range("D22").Value = "'0.910"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("D23").Value = "'18.30"
$ws.Range("E23").Value = "  +7.98%  "

$ws.Range("D24").Value = "'5.07"
$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("D25").Value = "'100.68"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("D26").Value = "'3.99"
$ws.Range("E26").Value = "  -2.23%  "

$ws.Range("D27").Value = "'5.99"
$ws.Range("E27").Value = "  -0.53%  "

$ws.Range("D28").Value = "'2.74"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("D29").Value = "'9.71"
$ws.Range("E29").Value = "  +3.11%  "

$ws.Range("D30").Value = "'8.62"
$ws.Range("E30").Value = "  -1.95%  "

$ws.Range("D31").Value = "'30.20"
$ws.Range("E31").Value = "  -1.93%  "

$ws.Range("D32").Value = "'6.70"
$ws.Range("E32").Value = "  +6.03%  "

$ws.Range("D33").Value = "'4.03"
$ws.Range("E33").Value = "  +4.50%  "

$ws.Range("D34").Value = "'563.78"
$ws.Range("E34").Value = "  +2.32%  "

$ws.Range("D35").Value = "'11.07"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("D36").Value = "'0.105"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'57.14"
$ws.Range("E37").Value = "  -1.56%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "'3.701.03"
$ws.Range("E39").Value = "  -2.68%  "

$ws.Range("D40").Value = "'0.0₃0725"
$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'33.99"
$ws.Range("E41").Value = "  +5.01%  "

$ws.Range("D42").Value = "'3.31"
$ws.Range("E42").Value = "  -3.14%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.130"
$ws.Range("E43").Value = "  +1.47%  "

$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.70"
$ws.Range("E44").Value = "  +0.54%  "

$ws.Range("B45").Value = "CoreDAO"
$ws.Range("C45").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D45").Value = "'3.35"
$ws.Range("E45").Value = "  -2.54%  "

$ws.Range("D46").Value = "'0.340"
$ws.Range("E46").Value = "  -2.63%  "

$ws.Range("D47").Value = "'0.0423"
$ws.Range("E47").Value = "  +1.02%  "

$ws.Range("D48").Value = "'3.24"
$ws.Range("E48").Value = "  +2.56%  "

$ws.Range("D49").Value = "'0.129"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").Value = "'2.59"
$ws.Range("E50").Value = "  -2.78%  "

$ws.Range("E51").Value = "  +0.07%  "
